$wb = $excel.ActiveWorkbook

# Sheet "zh-cn": row 16 holds the record for
# dff6ead6-...5dae3c21...zh-cn.xlf
# D16 = Correspond Handoff Datetime, G16 = Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D16").Value = "2016-03-04 11:12:45"
$wsZhCn.Range("G16").Value = "2016-03-04 11:13:49"

# Sheet "de-de": row 5 holds the record for
# dff6ead6-...5dae3c21...de-de.xlf
# D5 = Correspond Handoff Datetime, G5 = Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-04 11:13:01"
$wsDeDe.Range("G5").Value = "2016-03-04 11:14:23"
